# Add a new "2022" data column (column M) to the preschool-education
# coverage table, mirroring the formatting already used for the
# neighbouring 2020/2021 columns (K and L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> value for the new column M (row 4 is the year header, 5-30 are data)
$newColumnData = [ordered]@{
    4  = 2022
    5  = 24.6
    6  = 40.7
    7  = 20.7
    8  = 26.6
    9  = 44.5
    10 = 21.9
    11 = 21.9
    12 = 35.3
    13 = 17.6
    14 = 28
    15 = 44.9
    16 = 21.5
    17 = 36.2
    18 = 53.1
    19 = 33.4
    20 = 20.2
    21 = 15.4
    22 = 20.5
    23 = 27.1
    24 = 36.1
    25 = 25.2
    26 = 24.2
    27 = 46.5
    28 = 20.3
    29 = 40.5
    30 = 44.5
}

foreach ($row in $newColumnData.Keys) {
    $value = $newColumnData[$row]

    # Row 30 is the bottom "total" row, already bordered/number-formatted
    # via column L; every other row picks up the plain style already used
    # by column K in that row. Copying the format keeps fonts/borders/
    # number formats consistent with the rest of the table.
    if ($row -eq 30) {
        $ws.Range("L$row").Copy()
    } else {
        $ws.Range("K$row").Copy()
    }
    $ws.Range("M$row").PasteSpecial(-4122)
    $ws.Range("M$row").Value = $value
}

# Row 14 did not already have a one-decimal percentage style available on
# column K (it used a plain integer style), so give it the "0.0" number
# format explicitly - this creates the single new cellXfs entry.
$ws.Range("M14").NumberFormat = "0.0"

$excel.CutCopyMode = 0

# Match the author's final selection position recorded in the saved file.
$ws.Range("N7").Select()
